# Pre-Course Note.docx edit: add "Important Kali Linux Command;" heading
# after the "Kali linux ova file download and install:" bullet, and turn
# the trailing empty bookmarked list item into a bulleted (not numbered)
# paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that ends the "... ova file download and install: " line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("ova file download and install")) {
        $target = $p
    }
}

$afterRange = $target.Range
$afterRange.Collapse(0)

# Insert a blank paragraph right after it (matches the bare <w:p/> in the diff).
$blankPara = $afterRange.InsertParagraphAfter()

# Find the newly inserted blank paragraph and strip any inherited list
# numbering / style so it stays a plain paragraph.
$blankParaObj = $target.Next()
$blankParaObj.Range.ListFormat.RemoveNumbers()
$blankParaObj.Range.ParagraphFormat.Style = "Normal"

# Insert the bold "Important Kali Linux Command; " heading paragraph after that.
$headingInsertRange = $blankParaObj.Range
$headingInsertRange.Collapse(0)
$headingInsertRange.InsertParagraphAfter()

$headingPara = $blankParaObj.Next()
$headingPara.Range.Bold = 1
$headingPara.Range.BoldBi = 1
$headingPara.Range.Text = "Important Kali Linux Command; "

# Now find the empty bookmarked list paragraph ("_GoBack") that follows the
# "Kali linux ova file..." bullet and switch it from the numbered list
# (numId 4) onto a new bulleted list (numId 5), keeping the ListParagraph
# style.
$bookmarkPara = $headingPara.Next()
$bookmarkPara.Range.ListFormat.RemoveNumbers()

$gallery = $word.ListGalleries.Item(1)
$bulletTemplate = $gallery.ListTemplates.Item(1)
$bookmarkPara.Range.ListFormat.ApplyListTemplate($bulletTemplate)

Write-Output "edit complete"
